$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R data for year 2021
$ws.Range("R4").Value = 2021
$ws.Range("R4").Style = $ws.Range("Q4").Style

$ws.Range("R5").Value = 20.5
$ws.Range("R5").Style = $ws.Range("Q5").Style

# Update selection on the sheet view to match the new target cell
$ws.Range("S12").Select()
